$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill of Materials-BOM")

# J2 and J3 were previously blank; set them to match I2/I3 (a single-space text
# value). The leading apostrophe matches the text/quote-prefix formatting that
# I2/I3 already carry (cells entered as a plain space are stored as text).
$ws.Range("J2").Formula = "' "
$ws.Range("J3").Formula = "' "
